# Auto-generated edit script: applies market-price data refresh to Malboro_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# with values pulled from a scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1400493.6  # H17: 1433803.1 -> 1400493.6
$ws.Cells.Item(17, 10).Value = 1400493.6  # J17: 1433803.1 -> 1400493.6
$ws.Cells.Item(17, 12).Value = 4201480.800000001  # L17: 4301409.300000001 -> 4201480.800000001
$ws.Cells.Item(17, 14).Value = -4201816.800000001  # N17: -4301745.300000001 -> -4201816.800000001
$ws.Cells.Item(19, 8).Value = 1819.1333  # H19: 1828.5 -> 1819.1333
$ws.Cells.Item(19, 9).Value = 1805.7  # I19: 1818.7778 -> 1805.7
$ws.Cells.Item(19, 11).Value = 1805.7  # K19: 1818.7778 -> 1805.7
$ws.Cells.Item(19, 13).Value = -1630.7  # M19: -1643.7778 -> -1630.7
$ws.Cells.Item(62, 8).Value = 7900.3076  # H62: 7671.7144 -> 7900.3076
$ws.Cells.Item(62, 9).Value = 7080.625  # I62: 6816.1113 -> 7080.625
$ws.Cells.Item(62, 11).Value = 7080.625  # K62: 6816.1113 -> 7080.625
$ws.Cells.Item(62, 13).Value = -6456.625  # M62: -6192.1113 -> -6456.625
$ws.Cells.Item(65, 8).Value = 7900.3076  # H65: 7671.7144 -> 7900.3076
$ws.Cells.Item(65, 9).Value = 7080.625  # I65: 6816.1113 -> 7080.625
$ws.Cells.Item(65, 11).Value = 35403.125  # K65: 34080.5565 -> 35403.125
$ws.Cells.Item(65, 13).Value = -32283.125  # M65: -30960.5565 -> -32283.125
$ws.Cells.Item(69, 8).Value = 11350  # H69: 10233.167 -> 11350
$ws.Cells.Item(69, 10).Value = 12633.333  # J69: 10779.8 -> 12633.333
$ws.Cells.Item(69, 12).Value = 37899.999  # L69: 32339.4 -> 37899.999
$ws.Cells.Item(69, 14).Value = -39647.999  # N69: -34087.39999999999 -> -39647.999
$ws.Cells.Item(72, 8).Value = 11350  # H72: 10233.167 -> 11350
$ws.Cells.Item(72, 10).Value = 12633.333  # J72: 10779.8 -> 12633.333
$ws.Cells.Item(72, 12).Value = 113699.997  # L72: 97018.2 -> 113699.997
$ws.Cells.Item(72, 14).Value = -122435.997  # N72: -105754.2 -> -122435.997
$ws.Cells.Item(96, 8).Value = 1041.3334  # H96: 921.1429000000001 -> 1041.3334
$ws.Cells.Item(96, 10).Value = 500  # J96: 350 -> 500
$ws.Cells.Item(96, 12).Value = 1500  # L96: 1050 -> 1500
$ws.Cells.Item(96, 14).Value = -4246  # N96: -3796 -> -4246
$ws.Cells.Item(98, 8).Value = 15000  # H98: 2314.2 -> 15000
$ws.Cells.Item(98, 9).Value = 15000  # I98: 2383.3684 -> 15000
$ws.Cells.Item(98, 10).Value = 0  # J98: 1000 -> 0
$ws.Cells.Item(98, 11).Value = 15000  # K98: 2383.3684 -> 15000
$ws.Cells.Item(98, 12).Value = 0  # L98: 1000 -> 0
$ws.Cells.Item(98, 13).Value = -13502  # M98: -885.3683999999998 -> -13502
$ws.Cells.Item(98, 14).ClearContents()  # N98: was -3996
$ws.Cells.Item(122, 8).Value = 15000  # H122: 2314.2 -> 15000
$ws.Cells.Item(122, 9).Value = 15000  # I122: 2383.3684 -> 15000
$ws.Cells.Item(122, 10).Value = 0  # J122: 1000 -> 0
$ws.Cells.Item(122, 11).Value = 45000  # K122: 7150.1052 -> 45000
$ws.Cells.Item(122, 12).Value = 0  # L122: 3000 -> 0
$ws.Cells.Item(122, 13).Value = -42550  # M122: -4700.1052 -> -42550
$ws.Cells.Item(122, 14).ClearContents()  # N122: was -7900
$ws.Cells.Item(123, 8).Value = 114155.8  # H123: 142694.75 -> 114155.8
$ws.Cells.Item(123, 10).Value = 114155.8  # J123: 142694.75 -> 114155.8
$ws.Cells.Item(123, 12).Value = 114155.8  # L123: 142694.75 -> 114155.8
$ws.Cells.Item(123, 14).Value = -123955.8  # N123: -152494.75 -> -123955.8
$ws.Cells.Item(137, 8).Value = 8097.4375  # H137: 7930.8774 -> 8097.4375
$ws.Cells.Item(137, 9).Value = 1279.4584  # I137: 1225.72 -> 1279.4584
$ws.Cells.Item(137, 11).Value = 3838.3752  # K137: 3677.16 -> 3838.3752
$ws.Cells.Item(137, 13).Value = -1288.3752  # M137: -1127.16 -> -1288.3752
$ws.Cells.Item(141, 8).Value = 3352.1428  # H141: 1897.7 -> 3352.1428
$ws.Cells.Item(141, 9).Value = 2113.2  # I141: 1342 -> 2113.2
$ws.Cells.Item(141, 10).Value = 6449.5  # J141: 6899 -> 6449.5
$ws.Cells.Item(141, 11).Value = 6339.599999999999  # K141: 4026 -> 6339.599999999999
$ws.Cells.Item(141, 12).Value = 19348.5  # L141: 20697 -> 19348.5
$ws.Cells.Item(141, 13).Value = -1159.599999999999  # M141: 1154 -> -1159.599999999999
$ws.Cells.Item(141, 14).Value = -29708.5  # N141: -31057 -> -29708.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6304.0464  # H32: 6352.744 -> 6304.0464
$ws.Cells.Item(32, 9).Value = 5049.375  # I32: 5101.725 -> 5049.375
$ws.Cells.Item(32, 11).Value = 5049.375  # K32: 5101.725 -> 5049.375
$ws.Cells.Item(32, 13).Value = -4762.375  # M32: -4814.725 -> -4762.375
$ws.Cells.Item(97, 8).Value = 1229.6364  # H97: 1354.9678 -> 1229.6364
$ws.Cells.Item(97, 9).Value = 1170.5385  # I97: 1291.4 -> 1170.5385
$ws.Cells.Item(97, 10).Value = 1449.1428  # J97: 1619.8334 -> 1449.1428
$ws.Cells.Item(97, 11).Value = 1170.5385  # K97: 1291.4 -> 1170.5385
$ws.Cells.Item(97, 12).Value = 1449.1428  # L97: 1619.8334 -> 1449.1428
$ws.Cells.Item(97, 13).Value = -674.5385000000001  # M97: -795.4000000000001 -> -674.5385000000001
$ws.Cells.Item(97, 14).Value = -2441.1428  # N97: -2611.8334 -> -2441.1428

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 13121.108  # H31: 12616 -> 13121.108
$ws.Cells.Item(31, 9).Value = 6258.6665  # I31: 5801.2607 -> 6258.6665
$ws.Cells.Item(31, 11).Value = 6258.6665  # K31: 5801.2607 -> 6258.6665
$ws.Cells.Item(31, 13).Value = -5963.6665  # M31: -5506.2607 -> -5963.6665
$ws.Cells.Item(34, 8).Value = 13121.108  # H34: 12616 -> 13121.108
$ws.Cells.Item(34, 9).Value = 6258.6665  # I34: 5801.2607 -> 6258.6665
$ws.Cells.Item(34, 11).Value = 6258.6665  # K34: 5801.2607 -> 6258.6665
$ws.Cells.Item(34, 13).Value = -6056.6665  # M34: -5599.2607 -> -6056.6665
$ws.Cells.Item(70, 8).Value = 16666.666  # H70: 17000 -> 16666.666
$ws.Cells.Item(70, 10).Value = 16666.666  # J70: 17000 -> 16666.666
$ws.Cells.Item(70, 12).Value = 16666.666  # L70: 17000 -> 16666.666
$ws.Cells.Item(70, 14).Value = -17296.666  # N70: -17630 -> -17296.666
$ws.Cells.Item(73, 8).Value = 16666.666  # H73: 17000 -> 16666.666
$ws.Cells.Item(73, 10).Value = 16666.666  # J73: 17000 -> 16666.666
$ws.Cells.Item(73, 12).Value = 16666.666  # L73: 17000 -> 16666.666
$ws.Cells.Item(73, 14).Value = -18850.666  # N73: -19184 -> -18850.666
$ws.Cells.Item(99, 8).Value = 7063.926  # H99: 7224.077 -> 7063.926
$ws.Cells.Item(99, 9).Value = 3343.8  # I99: 3393.111 -> 3343.8
$ws.Cells.Item(99, 11).Value = 3343.8  # K99: 3393.111 -> 3343.8
$ws.Cells.Item(99, 13).Value = -1845.8  # M99: -1895.111 -> -1845.8
$ws.Cells.Item(105, 8).Value = 13916.625  # H105: 8174.9287 -> 13916.625
$ws.Cells.Item(105, 9).Value = 33778  # I105: 12993.875 -> 33778
$ws.Cells.Item(105, 10).Value = 1999.8  # J105: 1749.6666 -> 1999.8
$ws.Cells.Item(105, 11).Value = 33778  # K105: 12993.875 -> 33778
$ws.Cells.Item(105, 12).Value = 1999.8  # L105: 1749.6666 -> 1999.8
$ws.Cells.Item(105, 13).Value = -32031  # M105: -11246.875 -> -32031
$ws.Cells.Item(105, 14).Value = -5493.8  # N105: -5243.6666 -> -5493.8
$ws.Cells.Item(126, 8).Value = 7063.926  # H126: 7224.077 -> 7063.926
$ws.Cells.Item(126, 9).Value = 3343.8  # I126: 3393.111 -> 3343.8
$ws.Cells.Item(126, 11).Value = 10031.4  # K126: 10179.333 -> 10031.4
$ws.Cells.Item(126, 13).Value = -7561.400000000001  # M126: -7709.332999999999 -> -7561.400000000001
$ws.Cells.Item(132, 8).Value = 9032.1  # H132: 10313.588 -> 9032.1
$ws.Cells.Item(132, 9).Value = 3415.4443  # I132: 4238 -> 3415.4443
$ws.Cells.Item(132, 11).Value = 10246.3329  # K132: 12714 -> 10246.3329
$ws.Cells.Item(132, 13).Value = -7716.332900000001  # M132: -10184 -> -7716.332900000001
$ws.Cells.Item(134, 8).Value = 41674984  # H134: 43486870 -> 41674984
$ws.Cells.Item(134, 9).Value = 2146.5833  # I134: 2204.2727 -> 2146.5833
$ws.Cells.Item(134, 11).Value = 6439.749899999999  # K134: 6612.8181 -> 6439.749899999999
$ws.Cells.Item(134, 13).Value = -3904.749899999999  # M134: -4077.8181 -> -3904.749899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1760.3684  # H34: 1722.85 -> 1760.3684
$ws.Cells.Item(34, 10).Value = 2203.5  # J34: 2095 -> 2203.5
$ws.Cells.Item(34, 12).Value = 6610.5  # L34: 6285 -> 6610.5
$ws.Cells.Item(34, 14).Value = -6778.5  # N34: -6453 -> -6778.5
$ws.Cells.Item(46, 8).Value = 250362.5  # H46: 300 -> 250362.5
$ws.Cells.Item(46, 9).Value = 250362.5  # I46: 325 -> 250362.5
$ws.Cells.Item(46, 10).Value = 0  # J46: 250 -> 0
$ws.Cells.Item(46, 11).Value = 751087.5  # K46: 975 -> 751087.5
$ws.Cells.Item(46, 12).Value = 0  # L46: 750 -> 0
$ws.Cells.Item(46, 13).Value = -750996.5  # M46: -884 -> -750996.5
$ws.Cells.Item(46, 14).ClearContents()  # N46: was -932
$ws.Cells.Item(68, 8).Value = 1449.9  # H68: 1363.091 -> 1449.9
$ws.Cells.Item(68, 10).Value = 1857  # J68: 1686.75 -> 1857
$ws.Cells.Item(68, 12).Value = 5571  # L68: 5060.25 -> 5571
$ws.Cells.Item(68, 14).Value = -7193  # N68: -6682.25 -> -7193
$ws.Cells.Item(71, 8).Value = 1449.9  # H71: 1363.091 -> 1449.9
$ws.Cells.Item(71, 10).Value = 1857  # J71: 1686.75 -> 1857
$ws.Cells.Item(71, 12).Value = 16713  # L71: 15180.75 -> 16713
$ws.Cells.Item(71, 14).Value = -24825  # N71: -23292.75 -> -24825
$ws.Cells.Item(107, 8).Value = 3908287  # H107: 3908287.5 -> 3908287
$ws.Cells.Item(107, 10).Value = 5210732.5  # J107: 5210733.5 -> 5210732.5
$ws.Cells.Item(107, 12).Value = 15632197.5  # L107: 15632200.5 -> 15632197.5
$ws.Cells.Item(107, 14).Value = -15636037.5  # N107: -15636040.5 -> -15636037.5
$ws.Cells.Item(137, 8).Value = 1605.375  # H137: 1554.2727 -> 1605.375
$ws.Cells.Item(137, 9).Value = 1150  # I137: 1080.8 -> 1150
$ws.Cells.Item(137, 10).Value = 1878.6  # J137: 1948.8334 -> 1878.6
$ws.Cells.Item(137, 11).Value = 3450  # K137: 3242.4 -> 3450
$ws.Cells.Item(137, 12).Value = 5635.799999999999  # L137: 5846.5002 -> 5635.799999999999
$ws.Cells.Item(137, 13).Value = 1650  # M137: 1857.6 -> 1650
$ws.Cells.Item(137, 14).Value = -15835.8  # N137: -16046.5002 -> -15835.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 4233  # H102: 3929.6667 -> 4233
$ws.Cells.Item(102, 9).Value = 4358.8  # I102: 4016.4546 -> 4358.8
$ws.Cells.Item(102, 11).Value = 4358.8  # K102: 4016.4546 -> 4358.8
$ws.Cells.Item(102, 13).Value = -2736.8  # M102: -2394.4546 -> -2736.8
$ws.Cells.Item(122, 8).Value = 3462.125  # H122: 2251.1333 -> 3462.125
$ws.Cells.Item(122, 9).Value = 3299.5  # I122: 2296.8 -> 3299.5
$ws.Cells.Item(122, 10).Value = 3950  # J122: 2159.8 -> 3950
$ws.Cells.Item(122, 11).Value = 9898.5  # K122: 6890.400000000001 -> 9898.5
$ws.Cells.Item(122, 12).Value = 11850  # L122: 6479.400000000001 -> 11850
$ws.Cells.Item(122, 13).Value = -7448.5  # M122: -4440.400000000001 -> -7448.5
$ws.Cells.Item(122, 14).Value = -16750  # N122: -11379.4 -> -16750
$ws.Cells.Item(132, 8).Value = 16103.3125  # H132: 11381.042 -> 16103.3125
$ws.Cells.Item(132, 9).Value = 11587.846  # I132: 8544.842000000001 -> 11587.846
$ws.Cells.Item(132, 10).Value = 35670.332  # J132: 22158.6 -> 35670.332
$ws.Cells.Item(132, 11).Value = 34763.538  # K132: 25634.526 -> 34763.538
$ws.Cells.Item(132, 12).Value = 107010.996  # L132: 66475.79999999999 -> 107010.996
$ws.Cells.Item(132, 13).Value = -32233.538  # M132: -23104.526 -> -32233.538
$ws.Cells.Item(132, 14).Value = -112070.996  # N132: -71535.79999999999 -> -112070.996

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 5427.2  # H40: 3788.8333 -> 5427.2
$ws.Cells.Item(40, 9).Value = 4801  # I40: 2166 -> 4801
$ws.Cells.Item(40, 10).Value = 5523.5386  # J40: 5029.8237 -> 5523.5386
$ws.Cells.Item(40, 11).Value = 4801  # K40: 2166 -> 4801
$ws.Cells.Item(40, 12).Value = 5523.5386  # L40: 5029.8237 -> 5523.5386
$ws.Cells.Item(40, 13).Value = -4665  # M40: -2030 -> -4665
$ws.Cells.Item(40, 14).Value = -5795.5386  # N40: -5301.8237 -> -5795.5386
$ws.Cells.Item(122, 8).Value = 5681.8  # H122: 5559.0386 -> 5681.8
$ws.Cells.Item(122, 10).Value = 7181.25  # J122: 6660 -> 7181.25
$ws.Cells.Item(122, 12).Value = 21543.75  # L122: 19980 -> 21543.75
$ws.Cells.Item(122, 14).Value = -26443.75  # N122: -24880 -> -26443.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 5687.6553  # H132: 5599.729 -> 5687.6553
$ws.Cells.Item(132, 9).Value = 3150.5366  # I132: 3087.4285 -> 3150.5366
$ws.Cells.Item(132, 11).Value = 9451.6098  # K132: 9262.2855 -> 9451.6098
$ws.Cells.Item(132, 13).Value = -6921.6098  # M132: -6732.2855 -> -6921.6098
